$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 540.3200000000001
$ws.Range("I17").Value = 355.5
$ws.Range("J17").Value = 575.5238000000001
$ws.Range("K17").Value = 1066.5
$ws.Range("L17").Value = 1726.5714
$ws.Range("M17").Value = -898.5
$ws.Range("N17").Value = -2062.5714
$ws.Range("H19").Value = 3219.182
$ws.Range("I19").Value = 5582.222
$ws.Range("K19").Value = 5582.222
$ws.Range("M19").Value = -5407.222
$ws.Range("H33").Value = 347.54166
$ws.Range("I33").Value = 361.34784
$ws.Range("K33").Value = 361.34784
$ws.Range("M33").Value = -132.34784
$ws.Range("H61").Value = 105.07143
$ws.Range("I61").Value = 105.07143
$ws.Range("K61").Value = 315.21429
$ws.Range("M61").Value = -143.21429
$ws.Range("H92").Value = 577.7368
$ws.Range("I92").Value = 379.3846
$ws.Range("J92").Value = 1007.5
$ws.Range("K92").Value = 379.3846
$ws.Range("L92").Value = 1007.5
$ws.Range("M92").Value = 868.6154
$ws.Range("N92").Value = -3503.5
$ws.Range("H96").Value = 2037.7576
$ws.Range("I96").Value = 1733.3125
$ws.Range("J96").Value = 2324.2942
$ws.Range("K96").Value = 5199.9375
$ws.Range("L96").Value = 6972.882599999999
$ws.Range("M96").Value = -3826.9375
$ws.Range("N96").Value = -9718.882599999999
$ws.Range("H116").Value = 810.52
$ws.Range("I116").Value = 763.1875
$ws.Range("J116").Value = 894.6667
$ws.Range("K116").Value = 763.1875
$ws.Range("L116").Value = 894.6667
$ws.Range("M116").Value = 2678.8125
$ws.Range("N116").Value = -7778.6667
$ws.Range("H129").Value = 14615.352
$ws.Range("I129").Value = 686.8
$ws.Range("J129").Value = 16791.688
$ws.Range("K129").Value = 2060.4
$ws.Range("L129").Value = 50375.064
$ws.Range("M129").Value = 2939.6
$ws.Range("N129").Value = -60375.064
$ws.Range("H137").Value = 1264.8334
$ws.Range("I137").Value = 1106.9546
$ws.Range("J137").Value = 3001.5
$ws.Range("K137").Value = 3320.8638
$ws.Range("L137").Value = 9004.5
$ws.Range("M137").Value = -770.8638000000001
$ws.Range("N137").Value = -14104.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15462.838
$ws.Range("I32").Value = 16016.392
$ws.Range("J32").Value = 8635.666999999999
$ws.Range("K32").Value = 16016.392
$ws.Range("L32").Value = 8635.666999999999
$ws.Range("M32").Value = -15729.392
$ws.Range("N32").Value = -9209.666999999999
$ws.Range("H61").Value = 1967.8125
$ws.Range("I61").Value = 940.5263
$ws.Range("K61").Value = 940.5263
$ws.Range("M61").Value = -728.5263
$ws.Range("H122").Value = 862
$ws.Range("I122").Value = 889
$ws.Range("J122").Value = 700
$ws.Range("K122").Value = 2667
$ws.Range("L122").Value = 2100
$ws.Range("M122").Value = -217
$ws.Range("N122").Value = -7000
$ws.Range("H132").Value = 3555.638
$ws.Range("I132").Value = 3833.5366
$ws.Range("J132").Value = 2885.4119
$ws.Range("K132").Value = 11500.6098
$ws.Range("L132").Value = 8656.235700000001
$ws.Range("M132").Value = -8970.6098
$ws.Range("N132").Value = -13716.2357
$ws.Range("H136").Value = 1967.8125
$ws.Range("I136").Value = 940.5263
$ws.Range("K136").Value = 2821.5789
$ws.Range("M136").Value = -271.5789
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 733.94446
$ws.Range("I99").Value = 533.3333
$ws.Range("K99").Value = 533.3333
$ws.Range("M99").Value = 964.6667
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9093787
$ws.Range("I31").Value = 2216.9375
$ws.Range("K31").Value = 2216.9375
$ws.Range("M31").Value = -1921.9375
$ws.Range("H34").Value = 9093787
$ws.Range("I34").Value = 2216.9375
$ws.Range("K34").Value = 2216.9375
$ws.Range("M34").Value = -2014.9375
$ws.Range("H105").Value = 1637.25
$ws.Range("I105").Value = 1259.8
$ws.Range("K105").Value = 1259.8
$ws.Range("M105").Value = 487.2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 957481.75
$ws.Range("I131").Value = 16090
$ws.Range("J131").Value = 1323578.5
$ws.Range("K131").Value = 48270
$ws.Range("L131").Value = 3970735.5
$ws.Range("M131").Value = -43230
$ws.Range("N131").Value = -3980815.5
$ws.Range("H132").Value = 1292.3667
$ws.Range("I132").Value = 947.2143
$ws.Range("K132").Value = 8524.9287
$ws.Range("M132").Value = -5994.9287
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 61.2
$ws.Range("I2").Value = 26
$ws.Range("J2").Value = 84.666664
$ws.Range("K2").Value = 26
$ws.Range("L2").Value = 84.666664
$ws.Range("M2").Value = 87
$ws.Range("N2").Value = -310.666664
$ws.Range("H5").Value = 5334
$ws.Range("I5").Value = 5334
$ws.Range("K5").Value = 5334
$ws.Range("M5").Value = -5222
$ws.Range("H102").Value = 1540
$ws.Range("I102").Value = 1356
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 1356
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = 266
$ws.Range("N102").Value = -5244
$ws.Range("H113").Value = 13159639
$ws.Range("I113").Value = 20835362
$ws.Range("J113").Value = 1256.4286
$ws.Range("K113").Value = 20835362
$ws.Range("L113").Value = 1256.4286
$ws.Range("M113").Value = -20833192
$ws.Range("N113").Value = -5596.4286
$ws.Range("H132").Value = 45004.914
$ws.Range("I132").Value = 50993.49
$ws.Range("J132").Value = 4083
$ws.Range("K132").Value = 152980.47
$ws.Range("L132").Value = 12249
$ws.Range("M132").Value = -150450.47
$ws.Range("N132").Value = -17309
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1885.3572
$ws.Range("I7").Value = 1150
$ws.Range("J7").Value = 2436.875
$ws.Range("K7").Value = 1150
$ws.Range("L7").Value = 2436.875
$ws.Range("M7").Value = -1038
$ws.Range("N7").Value = -2660.875
$ws.Range("H22").Value = 441.66666
$ws.Range("I22").Value = 300
$ws.Range("J22").Value = 583.3333
$ws.Range("K22").Value = 300
$ws.Range("L22").Value = 583.3333
$ws.Range("M22").Value = -5
$ws.Range("N22").Value = -1173.3333
$ws.Range("H27").Value = 441.66666
$ws.Range("I27").Value = 300
$ws.Range("J27").Value = 583.3333
$ws.Range("K27").Value = 300
$ws.Range("L27").Value = 583.3333
$ws.Range("M27").Value = -193
$ws.Range("N27").Value = -797.3333
$ws.Range("H122").Value = 17004
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 1885.3572
$ws.Range("I126").Value = 1150
$ws.Range("J126").Value = 2436.875
$ws.Range("K126").Value = 3450
$ws.Range("L126").Value = 7310.625
$ws.Range("M126").Value = -980
$ws.Range("N126").Value = -12250.625
$ws.Range("H133").Value = 16379.223
$ws.Range("J133").Value = 16379.223
$ws.Range("L133").Value = 16379.223
$ws.Range("N133").Value = -21439.223
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 206.25
$ws.Range("I107").Value = 211.81818
$ws.Range("K107").Value = 635.4545400000001
$ws.Range("M107").Value = 1284.54546
